# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column F corresponds to the "municipio-nombre" field. It used to be
# classified as a measure (medida) of type xsd:int; it is now re-curated
# as a dimension (dim) tied to sdmx-dimension:refArea, with its URI
# template value set to "URI-Municipio" (mirroring how provincia-nombre
# and comarca-nombre are already modeled).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("F3").Value = "dim"
$ws.Range("F4").Value = "URI-Municipio"
